{"js": "// Replace the two-digit-by-two-digit multiplication equations shown in the\n// answer table with a newly generated set of equations (text-for-text swap,\n// one cell at a time). Each old value is unique in the document, so a\n// literal `body.search()` hit is unambiguous and corresponds to exactly one\n// table cell.\n\nconst pairs = [\n  [\"27\u00d769=1863\", \"24\u00d774=1776\"],\n  [\"20\u00d797=1940\", \"36\u00d744=1584\"],\n  [\"78\u00d787=6786\", \"21\u00d763=1323\"],\n  [\"62\u00d767=4154\", \"45\u00d714=630\"],\n  [\"88\u00d719=1672\", \"49\u00d739=1911\"],\n  [\"80\u00d717=1360\", \"30\u00d757=1710\"],\n  [\"79\u00d774=5846\", \"64\u00d764=4096\"],\n  [\"21\u00d729=609\", \"59\u00d794=5546\"],\n  [\"21\u00d785=1785\", \"31\u00d745=1395\"],\n  [\"19\u00d733=627\", \"48\u00d711=528\"],\n  [\"35\u00d714=490\", \"68\u00d744=2992\"],\n  [\"51\u00d797=4947\", \"35\u00d721=735\"],\n  [\"40\u00d792=3680\", \"39\u00d789=3471\"],\n  [\"20\u00d729=580\", \"79\u00d747=3713\"],\n  [\"55\u00d734=1870\", \"17\u00d751=867\"],\n  [\"81\u00d712=972\", \"12\u00d760=720\"],\n  [\"50\u00d772=3600\", \"72\u00d758=4176\"],\n  [\"51\u00d725=1275\", \"99\u00d761=6039\"],\n  [\"75\u00d752=3900\", \"44\u00d761=2684\"],\n  [\"88\u00d726=2288\", \"64\u00d741=2624\"],\n  [\"53\u00d717=901\", \"40\u00d732=1280\"],\n  [\"27\u00d751=1377\", \"52\u00d789=4628\"],\n  [\"27\u00d799=2673\", \"46\u00d712=552\"],\n  [\"60\u00d711=660\", \"15\u00d742=630\"],\n  [\"79\u00d740=3160\", \"37\u00d765=2405\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-by-two-digit multiplication equations shown in the\n# answer table with a newly generated set of equations (text-for-text swap,\n# one cell at a time). Each old value is unique in the document, so a plain\n# Find/Replace (no wildcards) on the whole document body is safe and will\n# touch exactly one cell per pair.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"27\u00d769=1863\"; New = \"24\u00d774=1776\" },\n    @{ Old = \"20\u00d797=1940\"; New = \"36\u00d744=1584\" },\n    @{ Old = \"78\u00d787=6786\"; New = \"21\u00d763=1323\" },\n    @{ Old = \"62\u00d767=4154\"; New = \"45\u00d714=630\" },\n    @{ Old = \"88\u00d719=1672\"; New = \"49\u00d739=1911\" },\n    @{ Old = \"80\u00d717=1360\"; New = \"30\u00d757=1710\" },\n    @{ Old = \"79\u00d774=5846\"; New = \"64\u00d764=4096\" },\n    @{ Old = \"21\u00d729=609\";  New = \"59\u00d794=5546\" },\n    @{ Old = \"21\u00d785=1785\"; New = \"31\u00d745=1395\" },\n    @{ Old = \"19\u00d733=627\";  New = \"48\u00d711=528\" },\n    @{ Old = \"35\u00d714=490\";  New = \"68\u00d744=2992\" },\n    @{ Old = \"51\u00d797=4947\"; New = \"35\u00d721=735\" },\n    @{ Old = \"40\u00d792=3680\"; New = \"39\u00d789=3471\" },\n    @{ Old = \"20\u00d729=580\";  New = \"79\u00d747=3713\" },\n    @{ Old = \"55\u00d734=1870\"; New = \"17\u00d751=867\" },\n    @{ Old = \"81\u00d712=972\";  New = \"12\u00d760=720\" },\n    @{ Old = \"50\u00d772=3600\"; New = \"72\u00d758=4176\" },\n    @{ Old = \"51\u00d725=1275\"; New = \"99\u00d761=6039\" },\n    @{ Old = \"75\u00d752=3900\"; New = \"44\u00d761=2684\" },\n    @{ Old = \"88\u00d726=2288\"; New = \"64\u00d741=2624\" },\n    @{ Old = \"53\u00d717=901\";  New = \"40\u00d732=1280\" },\n    @{ Old = \"27\u00d751=1377\"; New = \"52\u00d789=4628\" },\n    @{ Old = \"27\u00d799=2673\"; New = \"46\u00d712=552\" },\n    @{ Old = \"60\u00d711=660\";  New = \"15\u00d742=630\" },\n    @{ Old = \"79\u00d740=3160\"; New = \"37\u00d765=2405\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
